# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 2.5
$ws.Range("G2").Value = 2.52
$ws.Range("H2").Value = 3.4
$ws.Range("K2").Value = 3.25
$ws.Range("O2").Value = 1.46
$ws.Range("S2").Value = 4.7

# Row 3
$ws.Range("G3").Value = 1.45
$ws.Range("J3").Value = 4.7
$ws.Range("AA3").Value = 570
$ws.Range("AI3").Value = 230
$ws.Range("AM3").Value = 320

# Row 4
$ws.Range("AB4").Value = 9.199999999999999
$ws.Range("AD4").Value = 21
$ws.Range("AN4").Value = 8.6
